$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first row was a merged title row ("p4h ") spanning A1:E1.
# Remove it entirely; this shifts every other row up by one, drops the
# now-unused merge, and the now-unreferenced shared string is dropped
# automatically on save.
$ws.Rows.Item(1).Delete()

# Rows 1, 2 and 7 (old rows 2, 3, 8) don't use a fixed/custom row height -
# Excel auto-sizes them from their wrapped text, so after the shift their
# rendered height changes to match their new content. Rows 3-6 keep their
# previous explicit custom height untouched.
$ws.Rows.Item(1).RowHeight = 31.5
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Rows.Item(7).RowHeight = 15.75

# Restore the selection/active-cell state recorded in the saved file.
$ws.Range("A1:E1").Select()
